# Regenerated-data refresh (output generated at 456a3b4):
# bumps the "want-to-go" counts (column F) -- and a couple of "lowest
# price" values (column G) -- on all four sheets to their latest scraped
# totals.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1092
$ws.Range("F3").Value = 4704
$ws.Range("F5").Value = 195
$ws.Range("F6").Value = 1876
$ws.Range("F7").Value = 94
$ws.Range("F8").Value = 772
$ws.Range("F11").Value = 419
$ws.Range("F12").Value = 1146
$ws.Range("F13").Value = 1595
$ws.Range("F14").Value = 831
$ws.Range("F15").Value = 1866
$ws.Range("F16").Value = 578
$ws.Range("F17").Value = 532
$ws.Range("F18").Value = 626
$ws.Range("F19").Value = 207
$ws.Range("F20").Value = 26
$ws.Range("F21").Value = 26
$ws.Range("F23").Value = 1201
$ws.Range("F24").Value = 615
$ws.Range("F25").Value = 2528
$ws.Range("F28").Value = 1606
$ws.Range("F30").Value = 498
$ws.Range("F31").Value = 529
$ws.Range("F33").Value = 4310

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 39
$ws.Range("F18").Value = 285
$ws.Range("G19").Value = 90
$ws.Range("G20").Value = 90
$ws.Range("F22").Value = 145
$ws.Range("F29").Value = 4
$ws.Range("F30").Value = 25

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F7").Value = 322

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 1092
$ws.Range("F9").Value = 4704
$ws.Range("F11").Value = 195
$ws.Range("F12").Value = 1876
$ws.Range("F13").Value = 95
$ws.Range("F14").Value = 772
$ws.Range("F19").Value = 419
$ws.Range("F20").Value = 1146
$ws.Range("F22").Value = 39
$ws.Range("F25").Value = 831
$ws.Range("F26").Value = 1866
$ws.Range("F27").Value = 578
$ws.Range("F28").Value = 532
$ws.Range("F29").Value = 626
$ws.Range("F31").Value = 26
$ws.Range("F34").Value = 285
$ws.Range("G35").Value = 90
$ws.Range("F37").Value = 1201
$ws.Range("F38").Value = 145
$ws.Range("F39").Value = 2528
$ws.Range("F45").Value = 1606
$ws.Range("F46").Value = 498
$ws.Range("F49").Value = 4310
